# Update "想去人数" (interest count) values that changed in the source data refresh.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 3935
$wsExhibit.Range("F10").Value = 16
$wsExhibit.Range("F11").Value = 117

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 3935
$wsAll.Range("F11").Value = 16
$wsAll.Range("F12").Value = 117
